# Update computed profit-model values in the Pandaemonium Profits sheets.
# (currentAveragePrice / NQ / HQ, LevePrice NQ/HQ, LeveProfit NQ/HQ columns)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 320
$ws.Range("I18").Value = 320
$ws.Range("K18").Value = 320
$ws.Range("M18").Value = -36

$ws.Range("H21").Value = 4833.3335
$ws.Range("I21").Value = 4833.3335
$ws.Range("K21").Value = 4833.3335
$ws.Range("M21").Value = -4365.3335

$ws.Range("H23").Value = 4833.3335
$ws.Range("I23").Value = 4833.3335
$ws.Range("K23").Value = 4833.3335
$ws.Range("M23").Value = -4599.3335

$ws.Range("H86").Value = 85715.25
$ws.Range("I86").Value = 334494.34
$ws.Range("J86").Value = 2788.889
$ws.Range("K86").Value = 334494.34
$ws.Range("L86").Value = 2788.889
$ws.Range("M86").Value = -333371.34
$ws.Range("N86").Value = -5034.889

$ws.Range("H89").Value = 85715.25
$ws.Range("I89").Value = 334494.34
$ws.Range("J89").Value = 2788.889
$ws.Range("K89").Value = 1672471.7
$ws.Range("L89").Value = 13944.445
$ws.Range("M89").Value = -1666855.7
$ws.Range("N89").Value = -25176.445

$ws.Range("H141").Value = 3088.158
$ws.Range("I141").Value = 2729.625
$ws.Range("J141").Value = 5000.3335
$ws.Range("K141").Value = 8188.875
$ws.Range("L141").Value = 15001.0005
$ws.Range("M141").Value = -3008.875
$ws.Range("N141").Value = -25361.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1561.1111
$ws.Range("I2").Value = 1553.3334
$ws.Range("J2").Value = 1600
$ws.Range("K2").Value = 1553.3334
$ws.Range("L2").Value = 1600
$ws.Range("M2").Value = -1440.3334
$ws.Range("N2").Value = -1826

$ws.Range("H32").Value = 6370.7456
$ws.Range("I32").Value = 4762.755
$ws.Range("K32").Value = 4762.755
$ws.Range("M32").Value = -4475.755

$ws.Range("H45").Value = 1960.7273
$ws.Range("I45").Value = 1824.2222
$ws.Range("J45").Value = 2575
$ws.Range("K45").Value = 1824.2222
$ws.Range("L45").Value = 2575
$ws.Range("M45").Value = -1447.2222
$ws.Range("N45").Value = -3329

$ws.Range("H74").Value = 1284.6666
$ws.Range("I74").Value = 1373.4584
$ws.Range("J74").Value = 574.3333
$ws.Range("K74").Value = 1373.4584
$ws.Range("L74").Value = 574.3333
$ws.Range("M74").Value = -499.4584
$ws.Range("N74").Value = -2322.3333

$ws.Range("H77").Value = 1284.6666
$ws.Range("I77").Value = 1373.4584
$ws.Range("J77").Value = 574.3333
$ws.Range("K77").Value = 6867.291999999999
$ws.Range("L77").Value = 2871.6665
$ws.Range("M77").Value = -2499.291999999999
$ws.Range("N77").Value = -11607.6665

$ws.Range("H116").Value = 1561.1111
$ws.Range("I116").Value = 1553.3334
$ws.Range("J116").Value = 1600
$ws.Range("K116").Value = 1553.3334
$ws.Range("L116").Value = 1600
$ws.Range("M116").Value = 740.6666
$ws.Range("N116").Value = -6188

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1561.1111
$ws.Range("I3").Value = 1553.3334
$ws.Range("J3").Value = 1600
$ws.Range("K3").Value = 1553.3334
$ws.Range("L3").Value = 1600
$ws.Range("M3").Value = -1439.3334
$ws.Range("N3").Value = -1828

$ws.Range("H107").Value = 1450.1875
$ws.Range("I107").Value = 1293.0714
$ws.Range("K107").Value = 1293.0714
$ws.Range("M107").Value = 626.9286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 19999
$ws.Range("J17").Value = 19999
$ws.Range("L17").Value = 19999
$ws.Range("N17").Value = -20347

$ws.Range("H31").Value = 2981.1614
$ws.Range("I31").Value = 2194.611
$ws.Range("K31").Value = 2194.611
$ws.Range("M31").Value = -1899.611

$ws.Range("H34").Value = 2981.1614
$ws.Range("I34").Value = 2194.611
$ws.Range("K34").Value = 2194.611
$ws.Range("M34").Value = -1992.611

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 6668.75
$ws.Range("I87").Value = 3001
$ws.Range("J87").Value = 7891.3335
$ws.Range("K87").Value = 9003
$ws.Range("L87").Value = 23674.0005
$ws.Range("M87").Value = -7755
$ws.Range("N87").Value = -26170.0005

$ws.Range("H90").Value = 6668.75
$ws.Range("I90").Value = 3001
$ws.Range("J90").Value = 7891.3335
$ws.Range("K90").Value = 27009
$ws.Range("L90").Value = 71022.0015
$ws.Range("M90").Value = -20769
$ws.Range("N90").Value = -83502.0015

$ws.Range("H92").Value = 637.8461
$ws.Range("I92").Value = 562.5
$ws.Range("J92").Value = 671.3333
$ws.Range("K92").Value = 1687.5
$ws.Range("L92").Value = 2013.9999
$ws.Range("M92").Value = -439.5
$ws.Range("N92").Value = -4509.9999

$ws.Range("H106").Value = 1911.8334
$ws.Range("J106").Value = 1911.8334
$ws.Range("L106").Value = 5735.5002
$ws.Range("N106").Value = -7627.5002

$ws.Range("H110").Value = 1833.1923
$ws.Range("I110").Value = 588.5
$ws.Range("J110").Value = 2059.5
$ws.Range("K110").Value = 1765.5
$ws.Range("L110").Value = 6178.5
$ws.Range("M110").Value = 2324.5
$ws.Range("N110").Value = -14358.5

$ws.Range("H131").Value = 16120.5
$ws.Range("I131").Value = 545.25
$ws.Range("J131").Value = 21537.979
$ws.Range("K131").Value = 1635.75
$ws.Range("L131").Value = 64613.937
$ws.Range("M131").Value = 3404.25
$ws.Range("N131").Value = -74693.93700000001

$ws.Range("H134").Value = 3624.9487
$ws.Range("I134").Value = 3420.5293
$ws.Range("J134").Value = 3782.9092
$ws.Range("K134").Value = 10261.5879
$ws.Range("L134").Value = 11348.7276
$ws.Range("M134").Value = -5191.5879
$ws.Range("N134").Value = -21488.7276

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5412.529
$ws.Range("I70").Value = 5162.769
$ws.Range("K70").Value = 5162.769
$ws.Range("M70").Value = -4892.769

$ws.Range("H73").Value = 5412.529
$ws.Range("I73").Value = 5162.769
$ws.Range("K73").Value = 5162.769
$ws.Range("M73").Value = -4226.769

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 6602.4
$ws.Range("J11").Value = 9666.666999999999
$ws.Range("L11").Value = 9666.666999999999
$ws.Range("N11").Value = -9946.666999999999

$ws.Range("H46").Value = 974.5
$ws.Range("I46").Value = 950
$ws.Range("K46").Value = 950
$ws.Range("M46").Value = -762

$ws.Range("H132").Value = 4828
$ws.Range("I132").Value = 4569.5454
$ws.Range("K132").Value = 13708.6362
$ws.Range("M132").Value = -11178.6362

$ws.Range("H136").Value = 5392.467
$ws.Range("I136").Value = 3052.7827
$ws.Range("J136").Value = 13080
$ws.Range("K136").Value = 9158.348100000001
$ws.Range("L136").Value = 39240
$ws.Range("M136").Value = -6608.348100000001
$ws.Range("N136").Value = -44340

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4121
$ws.Range("I62").Value = 4000
$ws.Range("J62").Value = 4165
$ws.Range("K62").Value = 4000
$ws.Range("L62").Value = 4165
$ws.Range("M62").Value = -3376
$ws.Range("N62").Value = -5413

$ws.Range("H65").Value = 4121
$ws.Range("I65").Value = 4000
$ws.Range("J65").Value = 4165
$ws.Range("K65").Value = 20000
$ws.Range("L65").Value = 20825
$ws.Range("M65").Value = -16880
$ws.Range("N65").Value = -27065

$ws.Range("H132").Value = 2005
$ws.Range("I132").Value = 1244.7142
$ws.Range("J132").Value = 2596.3333
$ws.Range("K132").Value = 3734.1426
$ws.Range("L132").Value = 7788.999899999999
$ws.Range("M132").Value = -1204.1426
$ws.Range("N132").Value = -12848.9999

$ws.Range("H133").Value = 43307.5
$ws.Range("J133").Value = 43307.5
$ws.Range("L133").Value = 43307.5
$ws.Range("N133").Value = -53427.5
